$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)

$ws.Range("D2").Value = "39.521.94"
$ws.Range("E2").Value = "  +2.05%  "

$ws.Range("D3").Value = "2.161.75"
$ws.Range("E3").Value = "  +3.03%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'228.28"
$ws.Range("E5").Value = "  +0.30%  "

$ws.Range("D6").Value = "'0.634"
$ws.Range("E6").Value = "  +3.09%  "

$ws.Range("D7").Value = "'63.52"
$ws.Range("E7").Value = "  +2.34%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.394"
$ws.Range("E9").Value = "  +1.17%  "

$ws.Range("D10").Value = "'0.0852"
$ws.Range("E10").Value = "  +1.62%  "

$ws.Range("E11").Value = "  +0.54%  "

$ws.Range("D12").Value = "'15.99"
$ws.Range("E12").Value = "  +2.11%  "

$ws.Range("D13").Value = "2.482.30"
$ws.Range("E13").Value = "  +2.97%  "

$ws.Range("D14").Value = "'22.04"
$ws.Range("E14").Value = "  +0.07%  "

$ws.Range("E15").Value = "  -0.13%  "

$ws.Range("D16").Value = "'5.50"
$ws.Range("E16").Value = "  -0.62%  "

$ws.Range("D17").Value = "2.155.49"
$ws.Range("E17").Value = "  +0.62%  "

$ws.Range("D18").Value = "39.618.44"
$ws.Range("E18").Value = "  +2.37%  "

$ws.Range("D19").Value = "'6.19"
$ws.Range("E19").Value = "  +1.21%  "

$ws.Range("D20").Value = "'72.19"
$ws.Range("E20").Value = "  +0.84%  "

$ws.Range("D21").Value = "0.0₃0847"
$ws.Range("E21").Value = "  +1.03%  "

$ws.Range("D22").Value = "'229.71"
$ws.Range("E22").Value = "  +0.83%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "'2.38"
$ws.Range("E24").Value = "  +1.50%  "

$ws.Range("D25").Value = "'2.35"
$ws.Range("E25").Value = "  +1.69%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'172.04"
$ws.Range("E26").Value = "  +0.31%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'9.51"
$ws.Range("E27").Value = "  -1.12%  "

$ws.Range("D28").Value = "'0.139"
$ws.Range("E28").Value = "  +0.45%  "

$ws.Range("D29").Value = "'19.86"
$ws.Range("E29").Value = "  +2.94%  "

$ws.Range("E30").Value = "  +1.05%  "

$ws.Range("E31").Value = "  +4.83%  "

$ws.Range("E32").Value = "  +2.37%  "

$ws.Range("E33").Value = "  +1.23%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.70"
$ws.Range("E34").Value = "  -0.79%  "

$ws.Range("B35").Value = "THORChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D35").Value = "'7.02"
$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").Value = "'0.0618"
$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").Value = "'2.43"
$ws.Range("E37").Value = "  +1.85%  "

$ws.Range("D38").Value = "'3.65"
$ws.Range("E38").Value = "  +3.11%  "

$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.25%  "

$ws.Range("D40").Value = "'102.45"
$ws.Range("E40").Value = "  -0.19%  "

$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'17.83"
$ws.Range("E41").Value = "  -1.63%  "

$ws.Range("E42").Value = "  +0.17%  "

$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").Value = "'4.61"
$ws.Range("E43").Value = "  +10.96%  "

$ws.Range("D44").Value = "1.524.87"
$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("E45").Value = "  +0.36%  "

$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "'1.10"
$ws.Range("E46").Value = "  +4.97%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0924"
$ws.Range("E47").Value = "  +1.61%  "

$ws.Range("E48").Value = "  -0.06%  "

$ws.Range("D49").Value = "'7.73"
$ws.Range("E49").Value = "  -1.09%  "

$ws.Range("D50").Value = "'3.00"
$ws.Range("E50").Value = "  +0.98%  "

$ws.Range("D51").Value = "2.367.14"
$ws.Range("E51").Value = "  +3.06%  "
